$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    5  = -2
    10 = 0
    11 = 4
    18 = -10
    19 = -2
    20 = -4
    23 = -5
    24 = -3
    27 = 4
    31 = -10
    32 = -4
    34 = -2
    35 = 0
    40 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
